# Merge at the input level almost working. #286
#
# The input table has three "key" columns (Organizacion/D, Pais/G,
# Cuando/I) whose repeated values should be merged across the rows that
# share the same value, the same way "D5:D6" (Organizacion = UNICEF) was
# already merged in the source file:
#   - D5:D6 ("UNICEF")            -> already merged; just re-home its style
#   - G5:G6 ("Colombia")          -> newly merged, G6's duplicate value cleared
#   - I4:I7 ("1 March 2015")      -> newly merged, I5:I7 were blank already
#
# Row 5 also gained two numbers that had been missing (Hombres/Mujeres).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in the previously-missing counts on row 5 ---------------------
$ws.Range("E5").Value = 105
$ws.Range("F5").Value = 110

# --- merge the "Pais" pair (G5:G6) and the "Cuando" column (I4:I7) ------
# (D5:D6 is merged already in the source workbook, so it is left alone)
$ws.Range("G5:G6").Merge()
$ws.Range("I4:I7").Merge()

# --- re-style the merged ranges ------------------------------------------
# D/G merges: keep left/general horizontal alignment, move to top vertical
$xlTop = -4160
$xlCenter = -4108
$ws.Range("D5:D6").VerticalAlignment = $xlTop
$ws.Range("G5:G6").VerticalAlignment = $xlTop

# I merge: centered horizontally, top vertically
$ws.Range("I4:I7").HorizontalAlignment = $xlCenter
$ws.Range("I4:I7").VerticalAlignment = $xlTop

# --- move the active selection the way the editing session left it ------
$ws.Range("F6").Select() | Out-Null
